$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E ("Inscritos") values as per the diff
$ws.Range("E2").Value = 33
$ws.Range("E7").Value = 6
$ws.Range("E15").Value = 123
$ws.Range("E17").Value = 81
$ws.Range("E18").Value = 77
$ws.Range("E23").Value = 3
$ws.Range("E24").Value = 17
$ws.Range("E25").Value = 14
$ws.Range("E36").Value = 63
$ws.Range("E39").Value = 16

# Row 41: E, F, H all change
$ws.Range("E41").Value = 22
$ws.Range("F41").Value = 11
$ws.Range("H41").Value = 11

# Row 43: E, F, H all change
$ws.Range("E43").Value = 16
$ws.Range("F43").Value = 10
$ws.Range("H43").Value = 10

# Row 47: E, F, H all change
$ws.Range("E47").Value = 43
$ws.Range("F47").Value = 27
$ws.Range("H47").Value = 27

$ws.Range("E48").Value = 19
$ws.Range("E61").Value = 19
$ws.Range("E64").Value = 25
$ws.Range("E71").Value = 21
$ws.Range("E75").Value = 10
$ws.Range("E78").Value = 22
$ws.Range("E80").Value = 17
$ws.Range("E81").Value = 8
$ws.Range("E82").Value = 5
$ws.Range("E83").Value = 7
$ws.Range("E87").Value = 9
